$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed values shared by all rows 110-130 in this block (columns A,B,C,E,F,G,H,I,J,K,L,R)
$fixedA = 10
$fixedB = 'Vega Modelo de Temuco'
$fixedC = 'La Araucanía'
$fixedE = 9
$fixedF = 'Fruta'
$fixedG = 100108
$fixedH = 'Tropicales y subtropicales'
$fixedI = 100108004
$fixedJ = 'Papaya'
$fixedK = 'Cultivar IV Región'
$fixedL = 'Primera'
$fixedR = 'Provincia del Elquí'

# Write row 130 (brand-new row appended at the bottom), filling in the fixed
# (non-shifting) columns that are constant across this whole data block
$ws.Range("A130").Value = $fixedA
$ws.Range("B130").Value = $fixedB
$ws.Range("C130").Value = $fixedC
$ws.Range("E130").Value = $fixedE
$ws.Range("F130").Value = $fixedF
$ws.Range("G130").Value = $fixedG
$ws.Range("H130").Value = $fixedH
$ws.Range("I130").Value = $fixedI
$ws.Range("J130").Value = $fixedJ
$ws.Range("K130").Value = $fixedK
$ws.Range("L130").Value = $fixedL
$ws.Range("R130").Value = $fixedR

# Give the new D130 date cell the same number format as the cell above it
$ws.Range("D130").NumberFormat = $ws.Range("D129").NumberFormat

# The weekly update shifts each record's price/volume data (columns D, M, N, O, P,
# Q, S, T) down into the following row (oldest observation appended as new row 130),
# while row 110 receives this week's brand-new observation.
$ws.Range("D110").Value = 45218
$ws.Range("M110").Value = 100
$ws.Range("N110").Value = 24000
$ws.Range("O110").Value = 24000
$ws.Range("P110").Value = 24000
$ws.Range("Q110").Value = '$/bandeja 10 kilos'
$ws.Range("S110").Value = 2400
$ws.Range("T110").Value = 10

$ws.Range("D111").Value = 45142
$ws.Range("M111").Value = 55
$ws.Range("N111").Value = 25000
$ws.Range("O111").Value = 25000
$ws.Range("P111").Value = 25000
$ws.Range("Q111").Value = '$/bandeja 10 kilos'
$ws.Range("S111").Value = 2500
$ws.Range("T111").Value = 10

$ws.Range("D112").Value = 44242
$ws.Range("M112").Value = 55
$ws.Range("N112").Value = 25000
$ws.Range("O112").Value = 25000
$ws.Range("P112").Value = 25000
$ws.Range("Q112").Value = '$/bandeja 10 kilos'
$ws.Range("S112").Value = 2500
$ws.Range("T112").Value = 10

$ws.Range("D113").Value = 44447
$ws.Range("M113").Value = 20
$ws.Range("N113").Value = 20000
$ws.Range("O113").Value = 20000
$ws.Range("P113").Value = 20000
$ws.Range("Q113").Value = '$/bandeja 10 kilos'
$ws.Range("S113").Value = 2000
$ws.Range("T113").Value = 10

$ws.Range("D114").Value = 44265
$ws.Range("M114").Value = 40
$ws.Range("N114").Value = 21000
$ws.Range("O114").Value = 21000
$ws.Range("P114").Value = 21000
$ws.Range("Q114").Value = '$/bandeja 10 kilos'
$ws.Range("S114").Value = 2100
$ws.Range("T114").Value = 10

$ws.Range("D115").Value = 44362
$ws.Range("M115").Value = 50
$ws.Range("N115").Value = 20000
$ws.Range("O115").Value = 21000
$ws.Range("P115").Value = 20600
$ws.Range("Q115").Value = '$/bandeja 10 kilos'
$ws.Range("S115").Value = 2060
$ws.Range("T115").Value = 10

$ws.Range("D116").Value = 45148
$ws.Range("M116").Value = 80
$ws.Range("N116").Value = 25000
$ws.Range("O116").Value = 25000
$ws.Range("P116").Value = 25000
$ws.Range("Q116").Value = '$/bandeja 10 kilos'
$ws.Range("S116").Value = 2500
$ws.Range("T116").Value = 10

$ws.Range("D117").Value = 44491
$ws.Range("M117").Value = 40
$ws.Range("N117").Value = 20000
$ws.Range("O117").Value = 20000
$ws.Range("P117").Value = 20000
$ws.Range("Q117").Value = '$/bandeja 10 kilos'
$ws.Range("S117").Value = 2000
$ws.Range("T117").Value = 10

$ws.Range("D118").Value = 44601
$ws.Range("M118").Value = 80
$ws.Range("N118").Value = 25000
$ws.Range("O118").Value = 25000
$ws.Range("P118").Value = 25000
$ws.Range("Q118").Value = '$/bandeja 10 kilos'
$ws.Range("S118").Value = 2500
$ws.Range("T118").Value = 10

$ws.Range("D119").Value = 45028
$ws.Range("M119").Value = 200
$ws.Range("N119").Value = 3500
$ws.Range("O119").Value = 3500
$ws.Range("P119").Value = 3500
$ws.Range("Q119").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S119").Value = 3500
$ws.Range("T119").Value = 1

$ws.Range("D120").Value = 44623
$ws.Range("M120").Value = 95
$ws.Range("N120").Value = 23000
$ws.Range("O120").Value = 23000
$ws.Range("P120").Value = 23000
$ws.Range("Q120").Value = '$/bandeja 10 kilos'
$ws.Range("S120").Value = 2300
$ws.Range("T120").Value = 10

$ws.Range("D121").Value = 44468
$ws.Range("M121").Value = 40
$ws.Range("N121").Value = 20000
$ws.Range("O121").Value = 20000
$ws.Range("P121").Value = 20000
$ws.Range("Q121").Value = '$/bandeja 10 kilos'
$ws.Range("S121").Value = 2000
$ws.Range("T121").Value = 10

$ws.Range("D122").Value = 44487
$ws.Range("M122").Value = 80
$ws.Range("N122").Value = 20000
$ws.Range("O122").Value = 20000
$ws.Range("P122").Value = 20000
$ws.Range("Q122").Value = '$/bandeja 10 kilos'
$ws.Range("S122").Value = 2000
$ws.Range("T122").Value = 10

$ws.Range("D123").Value = 45215
$ws.Range("M123").Value = 100
$ws.Range("N123").Value = 24000
$ws.Range("O123").Value = 24000
$ws.Range("P123").Value = 24000
$ws.Range("Q123").Value = '$/bandeja 10 kilos'
$ws.Range("S123").Value = 2400
$ws.Range("T123").Value = 10

$ws.Range("D124").Value = 44973
$ws.Range("M124").Value = 35
$ws.Range("N124").Value = 42000
$ws.Range("O124").Value = 42000
$ws.Range("P124").Value = 42000
$ws.Range("Q124").Value = '$/caja 15 kilos granel'
$ws.Range("S124").Value = 2800
$ws.Range("T124").Value = 15

$ws.Range("D125").Value = 44382
$ws.Range("M125").Value = 45
$ws.Range("N125").Value = 21000
$ws.Range("O125").Value = 21000
$ws.Range("P125").Value = 21000
$ws.Range("Q125").Value = '$/bandeja 10 kilos'
$ws.Range("S125").Value = 2100
$ws.Range("T125").Value = 10

$ws.Range("D126").Value = 45173
$ws.Range("M126").Value = 150
$ws.Range("N126").Value = 24000
$ws.Range("O126").Value = 24000
$ws.Range("P126").Value = 24000
$ws.Range("Q126").Value = '$/bandeja 10 kilos'
$ws.Range("S126").Value = 2400
$ws.Range("T126").Value = 10

$ws.Range("D127").Value = 45140
$ws.Range("M127").Value = 55
$ws.Range("N127").Value = 25000
$ws.Range("O127").Value = 25000
$ws.Range("P127").Value = 25000
$ws.Range("Q127").Value = '$/bandeja 10 kilos'
$ws.Range("S127").Value = 2500
$ws.Range("T127").Value = 10

$ws.Range("D128").Value = 44980
$ws.Range("M128").Value = 80
$ws.Range("N128").Value = 40000
$ws.Range("O128").Value = 40000
$ws.Range("P128").Value = 40000
$ws.Range("Q128").Value = '$/caja 15 kilos granel'
$ws.Range("S128").Value = 2667
$ws.Range("T128").Value = 15

$ws.Range("D129").Value = 44460
$ws.Range("M129").Value = 30
$ws.Range("N129").Value = 20000
$ws.Range("O129").Value = 20000
$ws.Range("P129").Value = 20000
$ws.Range("Q129").Value = '$/bandeja 10 kilos'
$ws.Range("S129").Value = 2000
$ws.Range("T129").Value = 10

$ws.Range("D130").Value = 44392
$ws.Range("M130").Value = 80
$ws.Range("N130").Value = 20000
$ws.Range("O130").Value = 20000
$ws.Range("P130").Value = 20000
$ws.Range("Q130").Value = '$/bandeja 10 kilos'
$ws.Range("S130").Value = 2000
$ws.Range("T130").Value = 10
